$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Map of (rowIndex, colIndex) -> new text, addressed directly by table
# position so the many repeated "NN÷N=NN, N" style strings can't collide.
$updates = @(
    @{ Row = 1;  Col = 1; Text = "95÷8=11, 7" },
    @{ Row = 1;  Col = 2; Text = "42÷2=21, 0" },
    @{ Row = 1;  Col = 3; Text = "56÷4=14, 0" },
    @{ Row = 1;  Col = 4; Text = "87÷7=12, 3" },
    @{ Row = 1;  Col = 5; Text = "25÷2=12, 1" },

    @{ Row = 5;  Col = 1; Text = "61÷8=7, 5" },
    @{ Row = 5;  Col = 2; Text = "38÷3=12, 2" },
    @{ Row = 5;  Col = 3; Text = "56÷3=18, 2" },
    @{ Row = 5;  Col = 4; Text = "47÷6=7, 5" },
    @{ Row = 5;  Col = 5; Text = "48÷2=24, 0" },

    @{ Row = 9;  Col = 1; Text = "61÷2=30, 1" },
    @{ Row = 9;  Col = 2; Text = "98÷8=12, 2" },
    @{ Row = 9;  Col = 3; Text = "19÷5=3, 4" },
    @{ Row = 9;  Col = 4; Text = "53÷8=6, 5" },
    @{ Row = 9;  Col = 5; Text = "43÷8=5, 3" },

    @{ Row = 13; Col = 1; Text = "42÷5=8, 2" },
    @{ Row = 13; Col = 2; Text = "19÷3=6, 1" },
    @{ Row = 13; Col = 3; Text = "15÷6=2, 3" },
    @{ Row = 13; Col = 4; Text = "81÷5=16, 1" },
    @{ Row = 13; Col = 5; Text = "57÷8=7, 1" },

    @{ Row = 17; Col = 1; Text = "77÷6=12, 5" },
    @{ Row = 17; Col = 2; Text = "12÷5=2, 2" },
    @{ Row = 17; Col = 3; Text = "42÷5=8, 2" },
    @{ Row = 17; Col = 4; Text = "54÷3=18, 0" },
    @{ Row = 17; Col = 5; Text = "35÷4=8, 3" }
)

foreach ($u in $updates) {
    $cell = $table.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters so only the
    # visible text is replaced, leaving the cell/paragraph structure intact.
    $rng.MoveEnd(12, -1)
    $rng.Text = $u.Text
}
